# BPA cohort 3 update
# Update the date format within the TimeDateDura column (G) from
# "Jan 21, 2026" to "21 Jan, 2026" and widen column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtractedScans")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*Jan 21, 2026*") {
        $cell.Value2 = $val -replace "Jan 21, 2026", "21 Jan, 2026"
    }
}

$ws.Columns.Item(7).ColumnWidth = 27.285714285714285
